# Add a new "12-10-2020" data column (column AA) to the COVID19 time-series
# sheet, following the pattern already used by column Z (the previous day).
#
# Column AA holds:
#   - Row 1: the date label "12-10-2020" stored as literal text (like the
#            rest of the header row, which stores its dates as text rather
#            than as real date serials).
#   - Rows 2-36: the new case-count numbers for each state/UT.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header cell AA1 -------------------------------------------------
$header = $ws.Cells.Item(1, 27)

# Force Text formatting *before* assigning the value so Excel doesn't
# auto-convert the "DD-MM-YYYY" looking string into a date serial - this
# mirrors how the rest of the header row (D1:Z1) stores its dates as text.
$header.NumberFormat = "@"
$header.Value = "12-10-2020"

# Match the look of the neighbouring header cells (bold, centered, boxed)
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous

# ---- Data cells AA2:AA36 ----------------------------------------------
$newValues = @{
    2  = 3764
    3  = 703208
    4  = 9232
    5  = 164582
    6  = 184541
    7  = 11787
    8  = 113771
    9  = 3059
    10 = 281869
    11 = 33203
    12 = 132173
    13 = 130003
    14 = 14471
    15 = 71845
    16 = 83571
    17 = 580054
    18 = 191798
    19 = 4037
    20 = 129019
    21 = 1266240
    22 = 10504
    23 = 5142
    24 = 2010
    25 = 5743
    26 = 227615
    27 = 26291
    28 = 110865
    29 = 135990
    30 = 2920
    31 = 602038
    32 = 187342
    33 = 24403
    34 = 46931
    35 = 390566
    36 = 258948
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 27).Value = $newValues[$row]
}
